$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: B (coin name), C (link), D (price), E (volume %)
# "Text" flag marks D values that look like plain numbers and therefore need
# the cell pre-formatted as Text so Excel keeps the original literal digits
# (e.g. trailing zeros like "1.00" or "57.01") instead of coercing to a number.
$updates = @(
    @{ Row = 2; D = "70.908.02" },
    @{ Row = 3; D = "3.540.42"; E = "  -0.59%  " },
    @{ Row = 4; E = "  +0.05%  " },
    @{ Row = 5; D = "616.58"; E = "  +0.83%  "; DText = $true },
    @{ Row = 6; D = "174.54"; E = "  +0.85%  "; DText = $true },
    @{ Row = 7; D = "3.536.25"; E = "  -0.60%  " },
    @{ Row = 8; D = "0.611"; E = "  -1.03%  "; DText = $true },
    @{ Row = 9; D = "1.00"; E = "  +0.02%  "; DText = $true },
    @{ Row = 10; E = "  +1.23%  " },
    @{ Row = 11; D = "7.22"; E = "  -4.90%  "; DText = $true },
    @{ Row = 12; D = "0.587"; E = "  +0.23%  "; DText = $true },
    @{ Row = 13; D = "46.72"; E = "  +0.25%  "; DText = $true },
    @{ Row = 14; E = "  -0.08%  " },
    @{ Row = 15; D = "4.111.52"; E = "  -0.75%  " },
    @{ Row = 16; D = "8.44"; E = "  +0.68%  "; DText = $true },
    @{ Row = 17; D = "610.84"; E = "  -0.29%  "; DText = $true },
    @{ Row = 18; D = "3.547.18"; E = "  -0.47%  " },
    @{ Row = 19; D = "70.966.56"; E = "  +0.55%  " },
    @{ Row = 20; E = "  +1.24%  " },
    @{ Row = 21; D = "17.82"; E = "  +2.38%  "; DText = $true },
    @{ Row = 22; D = "0.887"; E = "  +0.21%  "; DText = $true },
    @{ Row = 23; D = "9.04"; E = "  -4.21%  "; DText = $true },
    @{ Row = 24; D = "15.71"; E = "  -2.29%  "; DText = $true },
    @{ Row = 25; D = "98.45"; E = "  +1.49%  "; DText = $true },
    @{ Row = 26; E = "  -1.08%  " },
    @{ Row = 27; E = "  +0.04%  " },
    @{ Row = 28; D = "2.59"; E = "  -1.27%  "; DText = $true },
    @{ Row = 29; D = "33.93"; E = "  +1.39%  "; DText = $true },
    @{ Row = 30; D = "9.13"; E = "  +0.64%  "; DText = $true },
    @{ Row = 31; D = "3.06"; E = "  +0.26%  "; DText = $true },
    @{ Row = 32; D = "8.17"; E = "  -3.99%  "; DText = $true },
    @{ Row = 33; E = "  +0.28%  " },
    @{ Row = 34; D = "6.87"; E = "  -1.55%  "; DText = $true },
    @{ Row = 35; D = "634.18"; E = "  +9.97%  "; DText = $true },
    @{ Row = 36; E = "  -1.00%  " },
    @{ Row = 37; E = "  +0.49%  " },
    @{ Row = 38; E = "  -4.05%  " },
    @{ Row = 39; D = "0.0474"; E = "  -1.88%  "; DText = $true },
    @{ Row = 40; D = "57.01"; DText = $true },
    @{ Row = 41; E = "  +0.03%  " },
    @{ Row = 42; E = "  +2.12%  " },
    @{ Row = 43; D = "0.0₃0742"; E = "  +5.24%  " },
    @{ Row = 44; D = "3.367.63"; E = "  -0.56%  " },
    @{ Row = 45; D = "2.99"; E = "  +0.42%  "; DText = $true },
    @{ Row = 46; E = "  -2.01%  " },
    @{ Row = 47; E = "  -3.02%  " },
    @{ Row = 48; E = "  -1.30%  " },
    @{ Row = 49; E = "  +0.50%  " },
    @{ Row = 50; D = "133.23"; E = "  -0.51%  "; DText = $true },
    @{ Row = 51; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "1.00"; E = "  -0.01%  "; DText = $true }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.ContainsKey("DText")) { $cell.NumberFormat = "@" }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
